$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(469, "PAL000019", "215 Test Ave", "Test", "PA"),
    @(470, "PAL000049", "216 Test Ave", "Test", "PA"),
    @(471, "PAL000001", "217 Test Ave", "Test", "PA"),
    @(472, "PAL000048", "218 Test Ave", "Test", "PA"),
    @(473, "PAL000002", "219 Test Ave", "Test", "PA"),
    @(474, "PAL000509", "220 Test Ave", "Test", "PA"),
    @(475, "PAL000003", "221 Test Ave", "Test", "PA"),
    @(476, "PAL000502", "222 Test Ave", "Test", "PA"),
    @(477, "PAL000265", "223 Test Ave", "Test", "PA"),
    @(478, "PAL000004", "224 Test Ave", "Test", "PA"),
    @(479, "PAL000266", "225 Test Ave", "Test", "PA"),
    @(480, "PAL000501", "226 Test Ave", "Test", "PA"),
    @(481, "PAL000005", "227 Test Ave", "Test", "PA"),
    @(482, "PAL000007", "228 Test Ave", "Test", "PA"),
    @(483, "PAL000008", "229 Test Ave", "Test", "PA"),
    @(484, "PAL000009", "230 Test Ave", "Test", "PA"),
    @(485, "PAL000047", "231 Test Ave", "Test", "PA"),
    @(486, "PAL000042", "232 Test Ave", "Test", "PA"),
    @(487, "PAL000013", "233 Test Ave", "Test", "PA"),
    @(488, "PAL000015", "234 Test Ave", "Test", "PA"),
    @(489, "PAL000016", "235 Test Ave", "Test", "PA"),
    @(490, "PAL000050", "236 Test Ave", "Test", "PA"),
    @(491, "PAL000017", "237 Test Ave", "Test", "PA"),
    @(492, "PAL000018", "238 Test Ave", "Test", "PA"),
    @(493, "PAL000496", "239 Test Ave", "Test", "PA"),
    @(494, "PAL000021", "240 Test Ave", "Test", "PA"),
    @(495, "PAL000022", "241 Test Ave", "Test", "PA"),
    @(496, "PAL000224", "242 Test Ave", "Test", "PA"),
    @(497, "PAL000499", "243 Test Ave", "Test", "PA"),
    @(498, "PAL000024", "244 Test Ave", "Test", "PA"),
    @(499, "PAL000026", "245 Test Ave", "Test", "PA"),
    @(500, "PAL000031", "246 Test Ave", "Test", "PA"),
    @(501, "PAL000033", "247 Test Ave", "Test", "PA"),
    @(502, "PAL000035", "248 Test Ave", "Test", "PA"),
    @(503, "PAL000481", "249 Test Ave", "Test", "PA"),
    @(504, "PAL000301", "250 Test Ave", "Test", "PA"),
    @(505, "PAL000137", "251 Test Ave", "Test", "PA"),
    @(506, "PAL000479", "252 Test Ave", "Test", "PA"),
    @(507, "PAL000036", "253 Test Ave", "Test", "PA"),
    @(508, "PAL000494", "254 Test Ave", "Test", "PA"),
    @(509, "PAL000039", "255 Test Ave", "Test", "PA"),
    @(510, "PAL000439", "256 Test Ave", "Test", "PA"),
    @(511, "PAL000040", "257 Test Ave", "Test", "PA"),
    @(512, "PAL000041", "258 Test Ave", "Test", "PA"),
    @(513, "PAL000490", "259 Test Ave", "Test", "PA"),
    @(514, "PAL000043", "260 Test Ave", "Test", "PA"),
    @(515, "PAL000440", "261 Test Ave", "Test", "PA"),
    @(516, "PAL000252", "262 Test Ave", "Test", "PA"),
    @(517, "state_nm_gov_susana_martinez", "263 Test Ave", "Test", "NM"),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
    $ws.Cells.Item($r, 5).Value2 = 11111
    $ws.Cells.Item($r, 5).NumberFormat = "00000"
    $ws.Cells.Item($r, 6).Value2 = 1111
    $ws.Cells.Item($r, 6).NumberFormat = "0000"
}

$null = $ws.Range("E515").Select()
